$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'261.02"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.05%"
$ws.Range("E2").Style = "Normal"
$ws.Range("E3").Value = "'-0.65%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'4.719"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-0.55%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.06218"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'2.26%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'6.742"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'1.17%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.8500"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'0.47%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9161"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.1403"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-0.34%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.05038"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'1.29%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07074"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-0.42%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.03099"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-1.12%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.09054"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.26%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001535"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.06%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.0006167"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'1.33%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005980"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-1.65%"
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'-0.26%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.170"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'0.68%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'-1.01%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'-0.67%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1311"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'1.06%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'4.096"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'0.11%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04224"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.42%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'-1.54%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004073"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'4.12%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001201"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'0.04%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'4.15%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D40").Value = "'0.03951"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'1.89%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.1113"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-0.19%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.004139"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'0.23%"
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'0.13%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01360"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-16.80%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005161"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-2.94%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'0.03%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D48").Value = "'0.2574"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'90.25%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'0.03%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'0.03%"
$ws.Range("E50").Style = "Normal"
